$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'القاعدة (مترجمه)'
$ws.Range('B2').Value = 'معاداة الشيوعية ,إس إف إن (مترجمه),إسلاموية ,أصولية إسلامية ,الفصائل (مترجمه),معاداة الصهيونية ,قائمة قابلة للطي (مترجمه),مكان (مترجمه),معاداة الاستعمارية ,المشاعر المعادية للهند (مترجمه),التحالف السني الشيعي (مترجمه),ويلسكي سيولو (مترجمه),الجهادية (مترجمه),معاداة السامية ,الفصائل الإسلامية القومية (مترجمه),معاداة أمريكا (مترجمه),ديوبندية ,التحالف السني الشيعي (مترجمه),السلفية ,ضد المثليين (مترجمه),الوحدة الاسلامية (مترجمه),الأصولية الإسلامية معاداة أمريكا معاداة الشيوعية معاداة الإمبريالية المشاعر المعادية للهنود معاداة المثليين معاداة السامية معاداة الإمبريالية الغربية معاداة الصهيونية قائمة بسيطة (مترجمه),الإمبريالية المعادية للغرب (مترجمه),الجهادية الديوبندية (مترجمه),قائمة بسيطة (مترجمه),غالاغر (مترجمه),رايت (مترجمه),مصدر (توضيح) ,حيز غير مقطوع ,جوناراتنا (مترجمه),برغن ,بيتر ل (مترجمه),حرب دينية ,داخل العالم السري لأسامة بن لادن (مترجمه),صحيفة نيويورك فري برس (مترجمه),ص ص- (مترجمه),الدعوة الوهابية ,التيار القطبي ,عالم الحديث (مترجمه),تم نفيه رسميا (مترجمه),مقدمة (مترجمه),ص ص (مترجمه),قطبية الجهادية وحدة المسلمين (مترجمه),مقدمة (مترجمه),ص ص (مترجمه),قائمة بسيطة (مترجمه),التعددية الإسلامية (مترجمه)'
$ws.Range('C2').Value = '-'
$ws.Range('D2').Value = ''

$ws.Range('A3').Value = 'هاى (مترجمه)'
$ws.Range('B3').Value = 'القومية اليسارية (مترجمه),اشتراكية ثورية ,الماركسيةاللينينية (مترجمه),قومية أرمنية '
$ws.Range('C3').Value = 'من اليسار إلى أقصى اليسار (مترجمه)'
$ws.Range('D3').Value = ''

$ws.Range('A4').Value = 'حركة حماس '
$ws.Range('B4').Value = 'إس إف إن (مترجمه),أصولية إسلامية ,مادة (مترجمه),معاداة الصهيونية ,ص ص (مترجمه),– إس إف إن (مترجمه),- مادة (مترجمه),أوبل (مترجمه),وطنية فلسطينية ,سياسة عسكرية ,دالاكورا (مترجمه),جيلفين (مترجمه),إس إف إن (مترجمه),مادة (مترجمه),إس إف إن (مترجمه),دونينج (مترجمه),ليتفاك (مترجمه),شرطة أفقية ,– إس إف إن (مترجمه),شرطة أفقية ,ستيبانوفا (مترجمه),شيما (مترجمه),مادة (مترجمه),القومية الاسلامية (مترجمه),إسلاموية '
$ws.Range('C4').Value = '-'
$ws.Range('D4').Value = 'يحيى السنوار ,خليل الحية (مترجمه),ابو عمر حسن (مترجمه),يحيى السنوار ,اغتيال محمد ضيف الله عفن (مترجمه),اغتيال محمد ضيف (مترجمه)'

$ws.Range('A5').Value = 'القوات اللبنانية '
$ws.Range('B5').Value = 'محافظة ليبرالية ,قومية لبنانية ,ديمقراطية مسيحية '
$ws.Range('C5').Value = 'يمينية '
$ws.Range('D5').Value = 'سمير جعجع ,بشير الجميل ,نائب الرئيس '

$ws.Range('A6').Value = 'القوات اللبنانية '
$ws.Range('B6').Value = 'معاداة الشيوعية ,الفيدرالية في لبنان (مترجمه),قومية لبنانية ,سياسة محافظة ,قومية مسيحية ,موارنة ,معاداة العروبة (مترجمه)'
$ws.Range('C6').Value = 'من اليمين إلى أقصى اليمين (مترجمه)'
$ws.Range('D6').Value = ''

$ws.Range('A7').Value = 'الجبهة اللبنانية '
$ws.Range('B7').Value = 'معاداة العرب (مترجمه),معاداة الشيوعية ,الفصائل (مترجمه),قومية فينيقية ,معاداة القومية العربية (مترجمه),معاداة فلسطين (مترجمه),قومية لبنانية ,قومية مسيحية ,الفلانخية '
$ws.Range('C7').Value = '-'
$ws.Range('D7').Value = ''

$ws.Range('A8').Value = 'منظمة التحرير الفلسطينية '
$ws.Range('B8').Value = 'القومية العربية (مترجمه),ماركسية ,فصائل أغسطس (مترجمه),علمانية ,الفصائل (مترجمه),أوبل (مترجمه),معاداة الصهيونية ,حل الدولة الواحدة (مترجمه),أغسطس البعثية الماركسية (مترجمه),وطنية فلسطينية ,معاداة الاستعمارية ,بعثية ,حل الدولة الواحدة معاداة الصهيونية (مترجمه),أغسطس ,سي إن (مترجمه),اشتراكية عربية ,قومية عربية '
$ws.Range('C8').Value = 'يسارية '
$ws.Range('D8').Value = ''

$ws.Range('A9').Value = 'منظمة التحرير الفلسطينية '
$ws.Range('B9').Value = 'القومية العربية (مترجمه),ماركسية ,فصائل أغسطس (مترجمه),علمانية ,الفصائل (مترجمه),أوبل (مترجمه),معاداة الصهيونية ,حل الدولة الواحدة (مترجمه),أغسطس البعثية الماركسية (مترجمه),وطنية فلسطينية ,معاداة الاستعمارية ,بعثية ,حل الدولة الواحدة معاداة الصهيونية (مترجمه),أغسطس ,سي إن (مترجمه),اشتراكية عربية ,قومية عربية '
$ws.Range('C9').Value = 'يسارية '
$ws.Range('D9').Value = ''

$ws.Range('A10').Value = 'حزب الكتائب اللبنانية (مترجمه)'
$ws.Range('B10').Value = 'قومية مسيحية ,معاداة الشيوعية ,محافظة اجتماعية ,مارونية سياسية ,قومية لبنانية ,ديمقراطية مسيحية ,الفلانخية '
$ws.Range('C10').Value = 'يمينية ,أقصى اليمين (مترجمه),يمين الوسط (مترجمه)'
$ws.Range('D10').Value = 'سامي الجميل ,بيار الجميل '

$ws.Range('A11').Value = 'جيش لبنان الجنوبي '
$ws.Range('B11').Value = 'معاداة الشيوعية ,الوحدة الاسلامية المسيحية (مترجمه),علمانية ,الفصائل (مترجمه),صهيونية ,معاداة فلسطين (مترجمه),التعددية المذهبية (مترجمه),مارونية سياسية ,قومية لبنانية '
$ws.Range('C11').Value = '-'
$ws.Range('D11').Value = ''

$ws.Range('A12').Value = 'الحزب السوري القومي الاجتماعي '
$ws.Range('B12').Value = 'معاداة السامية ,قومية اقتصادية ,معاداة الشيوعية ,قومية يسارية ,أوبل (مترجمه),قائمة قابلة للطي (مترجمه),التحرر السوري (مترجمه),القومية السورية ,فاشية '
$ws.Range('C12').Value = 'و يمين (مترجمه),في الماضي (مترجمه),لقد أخطأ البعض في اعتبار الحزب أو عناصر أيديولوجيته أو أعضائه منتمين إلى اليسار السياسي. (مترجمه),في بعض الأحيان يتم تصنيفها على أنها يمينية متطرفة (مترجمه),efn مزامنة (مترجمه),و يمين (مترجمه),في بعض الأحيان يتم تصنيفها على أنها يمينية متطرفة (مترجمه)'
$ws.Range('D12').Value = 'داء الكلب بنات (مترجمه)'

# Remove rows 13-19, which are no longer part of the table
$ws.Range("A13:D19").EntireRow.Delete()
